# Update column F (dSF) values for the specified rows to match the
# re-pulled data / mean calculation described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12 = 4
    14 = -9
    16 = 9
    19 = -7
    20 = -1
    21 = -2
    22 = 4
    24 = -3
    25 = -2
    28 = 2
    36 = 6
    38 = 5
    39 = 1
    44 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
